$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Start aging" column (E): header label, unit row, first data value.
# Write E2 before E1 so the shared-string table gets "sec" (7) then
# "Start aging" (8), matching how the column was authored.
$ws.Range("E2").Value = "sec"
$ws.Range("E1").Value = "Start aging"
$ws.Range("E3").Value = 40

$ws.Columns.Item(5).ColumnWidth = 10

$ws.Range("E3").Select()
